$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 17, shifting rows 17:83 down to 18:84
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with data
$ws.Cells.Item(17, 1).Value = 8
$ws.Cells.Item(17, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44910
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = 100114007
$ws.Cells.Item(17, 7).Value = "Jengibre"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 460
$ws.Cells.Item(17, 11).Value = 14000
$ws.Cells.Item(17, 12).Value = 15000
$ws.Cells.Item(17, 13).Value = 14500
$ws.Cells.Item(17, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(17, 15).Value = "Perú"
$ws.Cells.Item(17, 16).Value = 1115
$ws.Cells.Item(17, 17).Value = 13
$ws.Cells.Item(17, 18).Value = "Hortaliza"
